$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0
$ws.Range("A3").Formula = "=6*60"
$ws.Range("A4").Formula = "=6*60+A3"
$ws.Range("A5:A16").Formula = "=6*60+A4"
$ws.Range("A17").Formula = "=132*60+A16"

$ws.Range("A2:A17").NumberFormat = "0.0"

$ws.Range("A18").Select() | Out-Null
